$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-sort the worker arrears detail rows (B16:G30): group by worker, periods ascending.
# Net cell content is unchanged overall -- only row order/groupings move -- per the NIT-9001384477 update.
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1002200467"
$ws.Range("D16").Value = "ALEXIS JOSE BARRIOS JIMENEZ"
$ws.Range("E16").Value = "1607"
$ws.Range("F16").Value = 27578
$ws.Range("G16").Value = 781242

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1002200467"
$ws.Range("D17").Value = "ALEXIS JOSE BARRIOS JIMENEZ"
$ws.Range("E17").Value = "1608"
$ws.Range("F17").Value = 27578
$ws.Range("G17").Value = 781242

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "73121855"
$ws.Range("D18").Value = "ABRAHAM BELEÑO HERRERA"
$ws.Range("E18").Value = "1609"
$ws.Range("F18").Value = 27578
$ws.Range("G18").Value = 689454

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1002200467"
$ws.Range("D19").Value = "ALEXIS JOSE BARRIOS JIMENEZ"
$ws.Range("E19").Value = "1609"
$ws.Range("F19").Value = 27578
$ws.Range("G19").Value = 781242

$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "77164938"
$ws.Range("D20").Value = "YONIS ENRIQUE BANQUEZ DAZA"
$ws.Range("E20").Value = "1611"
$ws.Range("F20").Value = 5516
$ws.Range("G20").Value = 781242

$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "33101145"
$ws.Range("D21").Value = "NAIFI CABARCAS ANGULO"
$ws.Range("E21").Value = "1612"
$ws.Range("F21").Value = 91200
$ws.Range("G21").Value = 3600000

$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "1001803950"
$ws.Range("D22").Value = "RAFAEL GUILLERMO ROMERO MURILLO"
$ws.Range("E22").Value = "1704"
$ws.Range("F22").Value = 29509
$ws.Range("G22").Value = 781242

$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "1001803950"
$ws.Range("D23").Value = "RAFAEL GUILLERMO ROMERO MURILLO"
$ws.Range("E23").Value = "1705"
$ws.Range("F23").Value = 29509
$ws.Range("G23").Value = 781242

$ws.Range("B24").Value = "CC"
$ws.Range("C24").Value = "1001803950"
$ws.Range("D24").Value = "RAFAEL GUILLERMO ROMERO MURILLO"
$ws.Range("E24").Value = "1706"
$ws.Range("F24").Value = 29509
$ws.Range("G24").Value = 781242

$ws.Range("B25").Value = "CC"
$ws.Range("C25").Value = "1001803950"
$ws.Range("D25").Value = "RAFAEL GUILLERMO ROMERO MURILLO"
$ws.Range("E25").Value = "1707"
$ws.Range("F25").Value = 29509
$ws.Range("G25").Value = 781242

$ws.Range("B26").Value = "CC"
$ws.Range("C26").Value = "85445826"
$ws.Range("D26").Value = "HUGO ARMANDO RODRIGUEZ ARIAS"
$ws.Range("E26").Value = "1809"
$ws.Range("F26").Value = 90666
$ws.Range("G26").Value = 4000000

$ws.Range("B27").Value = "CC"
$ws.Range("C27").Value = "85445826"
$ws.Range("D27").Value = "HUGO ARMANDO RODRIGUEZ ARIAS"
$ws.Range("E27").Value = "1810"
$ws.Range("F27").Value = 160000
$ws.Range("G27").Value = 4000000

$ws.Range("B28").Value = "CC"
$ws.Range("C28").Value = "85445826"
$ws.Range("D28").Value = "HUGO ARMANDO RODRIGUEZ ARIAS"
$ws.Range("E28").Value = "1811"
$ws.Range("F28").Value = 160000
$ws.Range("G28").Value = 4000000

$ws.Range("B29").Value = "CC"
$ws.Range("C29").Value = "85445826"
$ws.Range("D29").Value = "HUGO ARMANDO RODRIGUEZ ARIAS"
$ws.Range("E29").Value = "1812"
$ws.Range("F29").Value = 160000
$ws.Range("G29").Value = 4000000

$ws.Range("B30").Value = "CC"
$ws.Range("C30").Value = "85445826"
$ws.Range("D30").Value = "HUGO ARMANDO RODRIGUEZ ARIAS"
$ws.Range("E30").Value = "1901"
$ws.Range("F30").Value = 160000
$ws.Range("G30").Value = 4000000
